$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (MuSCs -> Shh -> Boc -> ECs) ---
$ws.Cells.Item(2,5).Value  = 2
$ws.Cells.Item(2,6).Value  = 0.6666666666666666
$ws.Cells.Item(2,7).Value  = 0.08785
$ws.Cells.Item(2,8).Value  = 0.26355
$ws.Cells.Item(2,13).Value = 2.689040333333333
$ws.Cells.Item(2,14).Value = 8.067121
$ws.Cells.Item(2,15).Value = 0.1682412044246168
$ws.Cells.Item(2,16).Value = 0.1682412044246168
$ws.Cells.Item(2,17).Value = 0.2362321932833333
$ws.Cells.Item(2,18).Value = 2.12608973955
$ws.Cells.Item(2,19).Value = 0.1682412044246168
$ws.Cells.Item(2,20).Value = 0.1682412044246168

# --- Update existing row 3 (MuSCs -> Shh -> Boc -> FAPs) ---
$ws.Cells.Item(3,5).Value  = 2
$ws.Cells.Item(3,6).Value  = 0.6666666666666666
$ws.Cells.Item(3,7).Value  = 0.08785
$ws.Cells.Item(3,8).Value  = 0.26355
$ws.Cells.Item(3,15).Value = 0.7323297768450604
$ws.Cells.Item(3,16).Value = 0.7323297768450604
$ws.Cells.Item(3,17).Value = 1.0282847771
$ws.Cells.Item(3,18).Value = 9.2545629939
$ws.Cells.Item(3,19).Value = 0.7323297768450604
$ws.Cells.Item(3,20).Value = 0.7323297768450604

# --- Update existing row 4 (MuSCs -> Shh -> Boc -> MuSCs) ---
$ws.Cells.Item(4,5).Value  = 2
$ws.Cells.Item(4,6).Value  = 0.6666666666666666
$ws.Cells.Item(4,7).Value  = 0.08785
$ws.Cells.Item(4,8).Value  = 0.26355
$ws.Cells.Item(4,13).Value = 1.514375
$ws.Cells.Item(4,14).Value = 4.543125
$ws.Cells.Item(4,15).Value = 0.09474765803706024
$ws.Cells.Item(4,16).Value = 0.09474765803706024
$ws.Cells.Item(4,17).Value = 0.13303784375
$ws.Cells.Item(4,18).Value = 1.19734059375
$ws.Cells.Item(4,19).Value = 0.09474765803706024
$ws.Cells.Item(4,20).Value = 0.09474765803706024

# --- Add new row 5 (MuSCs -> Shh -> Boc -> Neutrophils) ---
$ws.Cells.Item(5,1).Value  = "MuSCs"
$ws.Cells.Item(5,2).Value  = "Shh"
$ws.Cells.Item(5,3).Value  = "Boc"
$ws.Cells.Item(5,4).Value  = "Neutrophils"
$ws.Cells.Item(5,5).Value  = 2
$ws.Cells.Item(5,6).Value  = 0.6666666666666666
$ws.Cells.Item(5,7).Value  = 0.08785
$ws.Cells.Item(5,8).Value  = 0.26355
$ws.Cells.Item(5,9).Value  = 1
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.07482333333333334
$ws.Cells.Item(5,14).Value = 0.22447
$ws.Cells.Item(5,15).Value = 0.004681360693262658
$ws.Cells.Item(5,16).Value = 0.004681360693262658
$ws.Cells.Item(5,17).Value = 0.006573229833333334
$ws.Cells.Item(5,18).Value = 0.0591590685
$ws.Cells.Item(5,19).Value = 0.004681360693262658
$ws.Cells.Item(5,20).Value = 0.004681360693262658

Write-Host "edits applied"
